$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.161.73"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.97"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.85"
$ws.Range("E5").Value = "  -2.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4642"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3863"
$ws.Range("E8").Value = "  -0.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07856"
$ws.Range("E9").Value = "  -0.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9632"
$ws.Range("E10").Value = "  -1.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.07"
$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.96"
$ws.Range("E12").Value = "  +3.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.694"
$ws.Range("E13").Value = "  -1.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.873"
$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06913"
$ws.Range("E15").Value = "  +0.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.53"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009954"
$ws.Range("E18").Value = "  -0.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.71"
$ws.Range("E19").Value = "  -1.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.162.19"
$ws.Range("E21").Value = "  -1.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.303"
$ws.Range("E22").Value = "  -1.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.03"
$ws.Range("E23").Value = "  -1.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.103"
$ws.Range("E24").Value = "  -2.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.051.21"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.44"
$ws.Range("E26").Value = "  +0.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.17"
$ws.Range("E27").Value = "  -1.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.756"
$ws.Range("E28").Value = "  -4.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.970"
$ws.Range("E29").Value = "  -1.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.81"
$ws.Range("E30").Value = "  +1.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09266"
$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9301"
$ws.Range("E32").Value = "  -3.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.289"
$ws.Range("E33").Value = "  -1.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.326"
$ws.Range("E34").Value = "  -1.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.333"
$ws.Range("E35").Value = "  -3.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05809"
$ws.Range("E36").Value = "  -4.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02100"
$ws.Range("E37").Value = "  -4.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.148"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.767"
$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5588"
$ws.Range("E40").Value = "  -1.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.884"
$ws.Range("E41").Value = "  -2.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1763"
$ws.Range("E42").Value = "  -1.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07230"
$ws.Range("E43").Value = "  +1.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.61"
$ws.Range("E44").Value = "  -1.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5272"
$ws.Range("E45").Value = "  -1.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.134"
$ws.Range("E46").Value = "  -7.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.126"
$ws.Range("E47").Value = "  -12.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.835"
$ws.Range("E48").Value = "  -3.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.79"
$ws.Range("E49").Value = "  +0.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.006"
$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.023"
$ws.Range("E51").Value = "  +0.21%  "

